$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: "Datos actualizados a 14 de Agosto de 2020 a las 14:06" -> "Datos actualizados a 14 de Agosto de 2020 a las 15:23"
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Agosto de 2020 a las 15:23"

# Row 4: "Estados Unidos" -> "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 5418045
$ws.Cells.Item(4, 3).Value = 2379
$ws.Cells.Item(4, 4).Value = 2844250
$ws.Cells.Item(4, 5).Value = 2403349
$ws.Cells.Item(4, 7).Value = 31
$ws.Cells.Item(4, 8).Value = 170446

# Row 16: "Arabia Saudita" -> "Arabia Saudita"
$ws.Cells.Item(16, 2).Value = 295902
$ws.Cells.Item(16, 3).Value = 1383
$ws.Cells.Item(16, 4).Value = 262959
$ws.Cells.Item(16, 5).Value = 29605
$ws.Cells.Item(16, 7).Value = 35
$ws.Cells.Item(16, 8).Value = 3338

# Row 28: "Catar" -> "Catar"
$ws.Cells.Item(28, 2).Value = 114532
$ws.Cells.Item(28, 3).Value = 251
$ws.Cells.Item(28, 4).Value = 111258
$ws.Cells.Item(28, 5).Value = 3084

# Row 42: "Bielorrusia" -> "Bielorrusia"
$ws.Cells.Item(42, 2).Value = 69308
$ws.Cells.Item(42, 3).Value = 105
$ws.Cells.Item(42, 4).Value = 66452
$ws.Cells.Item(42, 5).Value = 2253
$ws.Cells.Item(42, 7).Value = 4
$ws.Cells.Item(42, 8).Value = 603

# Row 44: "Emiratos Arabes Unidos" -> "Emiratos Arabes Unidos"
$ws.Cells.Item(44, 2).Value = 63819
$ws.Cells.Item(44, 3).Value = 330
$ws.Cells.Item(44, 4).Value = 57473
$ws.Cells.Item(44, 5).Value = 5987
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 359

# Row 45: "Paises Bajos" -> "Paises Bajos"
$ws.Cells.Item(45, 2).Value = 61840
$ws.Cells.Item(45, 3).Value = 636
$ws.Cells.Item(45, 7).Value = 2
$ws.Cells.Item(45, 8).Value = 6167

# Row 54: "Ghana" -> "Ghana"
$ws.Cells.Item(54, 2).Value = 41847
$ws.Cells.Item(54, 3).Value = 122
$ws.Cells.Item(54, 4).Value = 39718
$ws.Cells.Item(54, 5).Value = 1906

# Row 62: "Uzbekistan" -> "Uzbekistan"
$ws.Cells.Item(62, 2).Value = 33821
$ws.Cells.Item(62, 3).Value = 498
$ws.Cells.Item(62, 4).Value = 27825
$ws.Cells.Item(62, 5).Value = 5776
$ws.Cells.Item(62, 7).Value = 4
$ws.Cells.Item(62, 8).Value = 220

# Row 64: "Moldavia" -> "Kenia"
$ws.Cells.Item(64, 1).Value = "Kenia"
$ws.Cells.Item(64, 2).Value = 29334
$ws.Cells.Item(64, 3).Value = 580
$ws.Cells.Item(64, 4).Value = 15298
$ws.Cells.Item(64, 5).Value = 13571
$ws.Cells.Item(64, 7).Value = 5
$ws.Cells.Item(64, 8).Value = 465

# Row 65: "Serbia" -> "Moldavia"
$ws.Cells.Item(65, 1).Value = "Moldavia"
$ws.Cells.Item(65, 2).Value = 29087
$ws.Cells.Item(65, 4).Value = 20276
$ws.Cells.Item(65, 5).Value = 7933
$ws.Cells.Item(65, 8).Value = 878

# Row 66: "Kenia" -> "Serbia"
$ws.Cells.Item(66, 1).Value = "Serbia"
$ws.Cells.Item(66, 2).Value = 28998
$ws.Cells.Item(66, 4).Value = 26117
$ws.Cells.Item(66, 5).Value = 2220
$ws.Cells.Item(66, 8).Value = 661

# Row 103: "Croacia" -> "Croacia"
$ws.Cells.Item(103, 2).Value = 6258
$ws.Cells.Item(103, 3).Value = 208
$ws.Cells.Item(103, 4).Value = 5134
$ws.Cells.Item(103, 5).Value = 961
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 163

# Row 111: "Hong Kong" -> "Hong Kong"
$ws.Cells.Item(111, 5).Value = 902
$ws.Cells.Item(111, 7).Value = 2
$ws.Cells.Item(111, 8).Value = 67

# Row 139: "Nueva Zelanda" -> "Gambia"
$ws.Cells.Item(139, 1).Value = "Gambia"
$ws.Cells.Item(139, 2).Value = 1623
$ws.Cells.Item(139, 3).Value = 67
$ws.Cells.Item(139, 4).Value = 304
$ws.Cells.Item(139, 5).Value = 1269
$ws.Cells.Item(139, 7).Value = 7
$ws.Cells.Item(139, 8).Value = 50

# Row 140: "Gambia" -> "Nueva Zelanda"
$ws.Cells.Item(140, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(140, 2).Value = 1602
$ws.Cells.Item(140, 3).Value = 13
$ws.Cells.Item(140, 4).Value = 1531
$ws.Cells.Item(140, 5).Value = 49
$ws.Cells.Item(140, 8).Value = 22

# Row 213: "Islas Malvinas" -> "Montserrat"
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1

# Row 214: "Montserrat" -> "Islas Malvinas"
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0
